# Merge the split "}m through ", "n", " ", "of the preceding" runs in the
# "{m,n}" regex-quantifier bullet on slide 3 into a single run reading
# "}m through n of the preceding", keeping the formatting (rPr) that the
# original "n" run carried (lang="en-US" dirty="0" smtClean="0").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

$full = $tr.Text
$needle = "{m,n}m through n of the preceding"
$paraStart0 = $full.IndexOf($needle)

# Offsets (0-based) of each original run within the paragraph text:
#   "{"                 -> 0  (len 1)
#   "m,n"               -> 1  (len 3)
#   "}m through "       -> 4  (len 11)
#   "n"                 -> 15 (len 1)
#   " "                 -> 16 (len 1)
#   "of the preceding"  -> 17 (len 16)

# 1) Overwrite the isolated "n" run's text with the fully merged text.
#    Doing it on this single-character run keeps that run's own rPr
#    (dirty="0" smtClean="0") attached to the resulting text.
$runNStart0 = $paraStart0 + 15
$runNLen = 1
$sub = $tr.Characters($runNStart0 + 1, $runNLen)
$sub.Text = "}m through n of the preceding"

# 2) Remove the now-duplicated "}m through " text that precedes it.
$dupBeforeStart0 = $paraStart0 + 4
$dupBeforeLen = ("}m through ").Length
$before = $tr.Characters($dupBeforeStart0 + 1, $dupBeforeLen)
$before.Text = ""

# 3) Remove the now-duplicated " of the preceding" text that follows.
#    After step 2 the merged text shifted left by $dupBeforeLen characters.
$mergedStart0 = $paraStart0 + 4
$mergedLen = ("}m through n of the preceding").Length
$dupAfterStart0 = $mergedStart0 + $mergedLen
$dupAfterLen = (" of the preceding").Length
$after = $tr.Characters($dupAfterStart0 + 1, $dupAfterLen)
$after.Text = ""
